$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add hyperlinks to column E for the rows that now have a YouTube link.
# The order below matches the order the relationship ids (rId1..rId10)
# appear in the target workbook.
$ws.Hyperlinks.Add($ws.Range("E5"),  "https://www.youtube.com/watch?v=x-F7B1QVEII")
$ws.Hyperlinks.Add($ws.Range("E6"),  "https://www.youtube.com/watch?v=2TyIKM9-oZo")
$ws.Hyperlinks.Add($ws.Range("E8"),  "https://www.youtube.com/watch?v=Qswzk4eGPJg ")
$ws.Hyperlinks.Add($ws.Range("E9"),  "https://www.youtube.com/watch?v=rMrgTCOLnNs")
$ws.Hyperlinks.Add($ws.Range("E10"), "https://www.youtube.com/watch?v=-5w8A2LK950")
$ws.Hyperlinks.Add($ws.Range("E13"), "https://www.youtube.com/watch?v=x53ZI770W_4")
$ws.Hyperlinks.Add($ws.Range("E7"),  "https://www.youtube.com/watch?v=x-F7B1QVEII")
$ws.Hyperlinks.Add($ws.Range("E14"), "https://www.youtube.com/watch?v=lYWUwy_3_yU")
$ws.Hyperlinks.Add($ws.Range("E11"), "https://www.youtube.com/watch?v=AmvVawNlZ4A")
$ws.Hyperlinks.Add($ws.Range("E15"), "https://www.youtube.com/watch?v=OIdL1xAMWVs")

# Restore the view to show row 15 (as in the edited workbook).
$ws.Range("E15").Select()
